$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Hide master shapes on this slide (adds showMasterSp="0" to <p:sld>)
$s.DisplayMasterShapes = 0

# Remove the "Rectangle 1" shape (a plain white cover rectangle near the
# top-left corner of the slide) that was deleted in the authored edit.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 1") {
        $shp.Delete()
    }
}
